# Formed the consolidated report
# Update the "Absent" (column H) values to reflect the consolidated
# attendance report: a student is Absent (H=1) when there was no
# "Real" attendance recorded that day (column E = 0); otherwise H=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
